# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 120 (pushing the existing
# rows 120-191 down to 121-192) for "Vega Monumental Concepción - Piña".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 120, shifting the rest of
# the table (rows 120-191) down by one.
$ws.Rows.Item(120).Insert()

# Populate the newly inserted row 120 with the new weekly record.
$ws.Cells.Item(120, 1).Value = 11
$ws.Cells.Item(120, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(120, 3).Value = "Bíobío"
$ws.Cells.Item(120, 4).Value = 44777
$ws.Cells.Item(120, 5).Value = 8
$ws.Cells.Item(120, 6).Value = "Fruta"
$ws.Cells.Item(120, 7).Value = 100108
$ws.Cells.Item(120, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(120, 9).Value = 100108005
$ws.Cells.Item(120, 10).Value = "Piña"
$ws.Cells.Item(120, 11).Value = "Caramelo"
$ws.Cells.Item(120, 12).Value = "Segunda"
$ws.Cells.Item(120, 13).Value = 100
$ws.Cells.Item(120, 14).Value = 19000
$ws.Cells.Item(120, 15).Value = 20000
$ws.Cells.Item(120, 16).Value = 19500
$ws.Cells.Item(120, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(120, 18).Value = "Ecuador"
$ws.Cells.Item(120, 19).Value = 1393
$ws.Cells.Item(120, 20).Value = 14
